$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lesign ex 1")
$ws.Delete()
